$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newRow = 83

$ws.Cells.Item($newRow, 1).Value = "2025-05"
$ws.Cells.Item($newRow, 2).Value = 7
$ws.Cells.Item($newRow, 3).Value = 158
$ws.Cells.Item($newRow, 4).Value = 4.430379746835443
